$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6696
$ws1.Range("F3").Value = 797
$ws1.Range("F4").Value = 1099
$ws1.Range("F5").Value = 131
$ws1.Range("F6").Value = 703
$ws1.Range("F8").Value = 13
$ws1.Range("F9").Value = 1084
$ws1.Range("F10").Value = 844
$ws1.Range("F12").Value = 1322
$ws1.Range("F13").Value = 38
$ws1.Range("F16").Value = 540
$ws1.Range("F20").Value = 1478
$ws1.Range("F21").Value = 715
$ws1.Range("F22").Value = 284
$ws1.Range("F23").Value = 450
$ws1.Range("F24").Value = 438
$ws1.Range("F27").Value = 1123
$ws1.Range("F29").Value = 2343
$ws1.Range("F31").Value = 1204
$ws1.Range("F34").Value = 3818
$ws1.Range("F36").Value = 694

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 743
$ws2.Range("F19").Value = 4120
$ws2.Range("F32").Value = 53

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1244
$ws3.Range("F8").Value = 940

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1244
$ws4.Range("F6").Value = 940
$ws4.Range("F9").Value = 6696
$ws4.Range("F12").Value = 797
$ws4.Range("F13").Value = 743
$ws4.Range("F14").Value = 131
$ws4.Range("F15").Value = 703
$ws4.Range("F16").Value = 1084
$ws4.Range("F17").Value = 844
$ws4.Range("F22").Value = 1322
$ws4.Range("F23").Value = 38
$ws4.Range("F25").Value = 540
$ws4.Range("F28").Value = 1478
$ws4.Range("F29").Value = 715
$ws4.Range("F30").Value = 450
$ws4.Range("F31").Value = 438
$ws4.Range("F34").Value = 1123
$ws4.Range("F38").Value = 2343
$ws4.Range("F42").Value = 53
$ws4.Range("F45").Value = 1204
$ws4.Range("F48").Value = 3818
$ws4.Range("F50").Value = 694
